# Apply the Tue Jun 25 09:31:14 UTC 2024 cryptos-list refresh (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while keeping it TEXT even when it looks numeric
# (e.g. "570.04", "1.00"), without leaving the cell in a quote-prefixed style.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '60.879.65'
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").Value = '3.360.15'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.12%  '
Set-TextValue "D5" '570.04'
$ws.Range("E5").Value = '  -0.10%  '
Set-TextValue "D6" '135.13'
$ws.Range("E6").Value = '  +8.19%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '3.360.24'
$ws.Range("E8").Value = '  +0.09%  '
Set-TextValue "D9" '0.475'
$ws.Range("E9").Value = '  +0.14%  '
Set-TextValue "D10" '7.58'
$ws.Range("E10").Value = '  +5.64%  '
$ws.Range("E11").Value = '  +3.83%  '
$ws.Range("E12").Value = '  +4.26%  '
$ws.Range("D13").Value = '3.937.64'
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("E14").Value = '  +2.20%  '
Set-TextValue "D15" '0.0000171'
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").Value = '3.357.68'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("E17").Value = '  +2.64%  '
$ws.Range("D18").Value = '60.981.87'
$ws.Range("E18").Value = '  -2.24%  '
Set-TextValue "D19" '13.94'
$ws.Range("E19").Value = '  +6.99%  '
Set-TextValue "D20" '5.78'
$ws.Range("E20").Value = '  +4.39%  '
Set-TextValue "D21" '9.36'
$ws.Range("E21").Value = '  +3.41%  '
Set-TextValue "D22" '374.02'
$ws.Range("E22").Value = '  +4.05%  '
Set-TextValue "D23" '0.572'
$ws.Range("E23").Value = '  +4.21%  '
$ws.Range("D24").Value = '3.495.54'
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  +0.18%  '
Set-TextValue "D26" '70.65'
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  +12.01%  '
Set-TextValue "D28" '1.64'
$ws.Range("E28").Value = '  +16.47%  '
Set-TextValue "D29" '7.69'
$ws.Range("E29").Value = '  +11.28%  '
Set-TextValue "D30" '1.00'
$ws.Range("E30").Value = '  -0.23%  '
Set-TextValue "D31" '8.07'
$ws.Range("E31").Value = '  +3.83%  '
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("E33").Value = '  +4.88%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").Value = '3.395.66'
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +4.36%  '
Set-TextValue "D37" '5.53'
$ws.Range("E37").Value = '  +4.48%  '
$ws.Range("E38").Value = '  +4.64%  '
$ws.Range("E39").Value = '  +4.52%  '
Set-TextValue "D40" '163.38'
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("E41").Value = '  +4.99%  '
Set-TextValue "D42" '1.00'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D43" '41.44'
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue "D44" '1.20'
$ws.Range("E44").Value = '  +11.67%  '
Set-TextValue "D45" '4.37'
$ws.Range("E45").Value = '  +4.79%  '
$ws.Range("E46").Value = '  -0.28%  '
$ws.Range("E47").Value = '  +5.30%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D48" '22.87'
$ws.Range("E48").Value = '  +2.85%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D49" '6.94'
$ws.Range("E49").Value = '  +6.16%  '
Set-TextValue "D50" '23.01'
$ws.Range("E50").Value = '  +13.88%  '
Set-TextValue "D51" '2.41'
$ws.Range("E51").Value = '  +14.70%  '
